$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 7000
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Row 43
$ws.Range("H43").Value = 958.7143
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 958.7143
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 958.7143
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1096.7143

# Row 62
$ws.Range("H62").Value = 4174.3335
$ws.Range("I62").Value = 3025
$ws.Range("K62").Value = 3025
$ws.Range("M62").Value = -2401

# Row 65
$ws.Range("H65").Value = 4174.3335
$ws.Range("I65").Value = 3025
$ws.Range("K65").Value = 15125
$ws.Range("M65").Value = -12005

# Row 76
$ws.Range("H76").Value = 5053650.5
$ws.Range("I76").Value = 3433.3333
$ws.Range("J76").Value = 11113911
$ws.Range("K76").Value = 3433.3333
$ws.Range("L76").Value = 11113911
$ws.Range("M76").Value = -3118.3333
$ws.Range("N76").Value = -11114541

# Row 79
$ws.Range("H79").Value = 5053650.5
$ws.Range("I79").Value = 3433.3333
$ws.Range("J79").Value = 11113911
$ws.Range("K79").Value = 3433.3333
$ws.Range("L79").Value = 11113911
$ws.Range("M79").Value = -2341.3333
$ws.Range("N79").Value = -11116095

# Row 95
$ws.Range("H95").Value = 33250
$ws.Range("J95").Value = 33250
$ws.Range("L95").Value = 33250
$ws.Range("N95").Value = -38742

# Row 98
$ws.Range("H98").Value = 1170.2778
$ws.Range("I98").Value = 1004.6429
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 1004.6429
$ws.Range("L98").Value = 1750
$ws.Range("M98").Value = 493.3570999999999
$ws.Range("N98").Value = -4746

# Row 100
$ws.Range("H100").Value = 2207
$ws.Range("I100").Value = 1440
$ws.Range("K100").Value = 1440
$ws.Range("M100").Value = -899

# Row 103
$ws.Range("H103").Value = 66667010
$ws.Range("I103").Value = 166666860
$ws.Range("J103").Value = 433
$ws.Range("K103").Value = 500000580
$ws.Range("L103").Value = 1299
$ws.Range("M103").Value = -499999994
$ws.Range("N103").Value = -2471

# Row 122
$ws.Range("H122").Value = 1170.2778
$ws.Range("I122").Value = 1004.6429
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 3013.9287
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -563.9287000000004
$ws.Range("N122").Value = -10150

# Row 129
$ws.Range("H129").Value = 159629.03
$ws.Range("J129").Value = 189701.08
$ws.Range("L129").Value = 569103.24
$ws.Range("N129").Value = -579103.24

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3224.0833
$ws.Range("I61").Value = 2586.125
$ws.Range("K61").Value = 2586.125
$ws.Range("M61").Value = -2374.125

# Row 63
$ws.Range("H63").Value = 3908207.5
$ws.Range("I63").Value = 2237.1428
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2237.1428
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -1551.1428
$ws.Range("N63").Value = -31251372

# Row 66
$ws.Range("H66").Value = 3908207.5
$ws.Range("I66").Value = 2237.1428
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 11185.714
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -7753.714
$ws.Range("N66").Value = -156256864

# Row 110
$ws.Range("H110").Value = 1014.3889
$ws.Range("I110").Value = 889.6923
$ws.Range("J110").Value = 1338.6
$ws.Range("K110").Value = 889.6923
$ws.Range("L110").Value = 1338.6
$ws.Range("M110").Value = 1155.3077
$ws.Range("N110").Value = -5428.6

# Row 122
$ws.Range("H122").Value = 1436.2333
$ws.Range("I122").Value = 1376.4231
$ws.Range("J122").Value = 1825
$ws.Range("K122").Value = 4129.2693
$ws.Range("L122").Value = 5475
$ws.Range("M122").Value = -1679.2693
$ws.Range("N122").Value = -10375

# Row 136
$ws.Range("H136").Value = 3224.0833
$ws.Range("I136").Value = 2586.125
$ws.Range("K136").Value = 7758.375
$ws.Range("M136").Value = -5208.375

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 525
$ws.Range("I5").Value = 366.66666
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 366.66666
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -253.66666
$ws.Range("N5").Value = -1226

# Row 105
$ws.Range("H105").Value = 2001814.4
$ws.Range("I105").Value = 1637
$ws.Range("K105").Value = 1637
$ws.Range("M105").Value = 110

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1597.2858
$ws.Range("I16").Value = 1440.2
$ws.Range("K16").Value = 1440.2
$ws.Range("M16").Value = -1153.2

# Row 62
$ws.Range("H62").Value = 6084.3335
$ws.Range("I62").Value = 5375
$ws.Range("K62").Value = 5375
$ws.Range("M62").Value = -4751

# Row 65
$ws.Range("H65").Value = 6084.3335
$ws.Range("I65").Value = 5375
$ws.Range("K65").Value = 26875
$ws.Range("M65").Value = -23755

# Row 70
$ws.Range("H70").Value = 15000
$ws.Range("J70").Value = 15000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15630

# Row 73
$ws.Range("H73").Value = 15000
$ws.Range("J73").Value = 15000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -17184

# Row 113
$ws.Range("H113").Value = 1597.2858
$ws.Range("I113").Value = 1440.2
$ws.Range("K113").Value = 1440.2
$ws.Range("M113").Value = 729.8

# Row 141
$ws.Range("H141").Value = 30853.389
$ws.Range("J141").Value = 32197.705
$ws.Range("L141").Value = 32197.705
$ws.Range("N141").Value = -42557.705

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 263.75
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 263.75
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 791.25
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -1015.25

# Row 12
$ws.Range("H12").Value = 97.454544
$ws.Range("J12").Value = 105.2
$ws.Range("L12").Value = 315.6
$ws.Range("N12").Value = -661.6

# Row 33
$ws.Range("H33").Value = 112.5
$ws.Range("I33").Value = 4.5
$ws.Range("J33").Value = 134.1
$ws.Range("K33").Value = 27
$ws.Range("L33").Value = 804.5999999999999
$ws.Range("M33").Value = 256
$ws.Range("N33").Value = -1370.6

# Row 131
$ws.Range("H131").Value = 693.66
$ws.Range("J131").Value = 720.17584
$ws.Range("L131").Value = 2160.52752
$ws.Range("N131").Value = -12240.52752

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 29400
$ws.Range("J57").Value = 29400
$ws.Range("L57").Value = 29400
$ws.Range("N57").Value = -31040

# Row 80
$ws.Range("H80").Value = 4033.6667
$ws.Range("I80").Value = 3600
$ws.Range("J80").Value = 4309.636
$ws.Range("K80").Value = 3600
$ws.Range("L80").Value = 4309.636
$ws.Range("M80").Value = -2602
$ws.Range("N80").Value = -6305.636

# Row 83
$ws.Range("H83").Value = 4033.6667
$ws.Range("I83").Value = 3600
$ws.Range("J83").Value = 4309.636
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 21548.18
$ws.Range("M83").Value = -13008
$ws.Range("N83").Value = -31532.18

# Row 135
$ws.Range("H135").Value = 46144.668
$ws.Range("J135").Value = 46144.668
$ws.Range("L135").Value = 46144.668
$ws.Range("N135").Value = -56284.668

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1818.8572
$ws.Range("I82").Value = 1869.8
$ws.Range("J82").Value = 1691.5
$ws.Range("K82").Value = 1869.8
$ws.Range("L82").Value = 1691.5
$ws.Range("M82").Value = -1508.8
$ws.Range("N82").Value = -2413.5

# Row 85
$ws.Range("H85").Value = 1818.8572
$ws.Range("I85").Value = 1869.8
$ws.Range("J85").Value = 1691.5
$ws.Range("K85").Value = 1869.8
$ws.Range("L85").Value = 1691.5
$ws.Range("M85").Value = -621.8
$ws.Range("N85").Value = -4187.5

# Row 93
$ws.Range("H93").Value = 1526.6364
$ws.Range("I93").Value = 1698.8334
$ws.Range("K93").Value = 1698.8334
$ws.Range("M93").Value = -450.8334

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4571.857
$ws.Range("J62").Value = 4571.857
$ws.Range("L62").Value = 4571.857
$ws.Range("N62").Value = -5819.857

# Row 65
$ws.Range("H65").Value = 4571.857
$ws.Range("J65").Value = 4571.857
$ws.Range("L65").Value = 22859.285
$ws.Range("N65").Value = -29099.285

# Row 81
$ws.Range("H81").Value = 2811.9
$ws.Range("I81").Value = 2100.5
$ws.Range("J81").Value = 2989.75
$ws.Range("K81").Value = 4201
$ws.Range("L81").Value = 5979.5
$ws.Range("M81").Value = -3140
$ws.Range("N81").Value = -8101.5

# Row 84
$ws.Range("H84").Value = 2811.9
$ws.Range("I84").Value = 2100.5
$ws.Range("J84").Value = 2989.75
$ws.Range("K84").Value = 21005
$ws.Range("L84").Value = 29897.5
$ws.Range("M84").Value = -15701
$ws.Range("N84").Value = -40505.5
